$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Molecular Weight" column (D) with header + a constant 208 value
# for every data row (rows 2-15), matching the existing table layout.
$ws.Range("D1").Value = "Molecular Weight"
$ws.Range("D2:D15").Value = 208

# Reuse the existing data-cell formatting (font Arial) instead of creating
# a brand new style entry, by copying the format from an existing data
# cell and pasting only the formatting onto the new column.
$ws.Range("C2").Copy()
$ws.Range("D1:D15").PasteSpecial(-4122)

# Match the active selection left behind in the saved workbook.
$ws.Range("D4:D15").Select()
